# Word COM-interop script implementing the diff:
#  1) Title: "INDEFINIDO " -> "[TERMINO]"          (paragraph 1)
#  2) Table cell: "MANIZALES" -> "[CD_CONT]"        (place-of-work city)
#     Table cell: "CALDAS"    -> "[DPTO_CONT]"      (place-of-work department)
#  3) Clause body: "indefinido " -> "[TERMINO]" + a
#     freshly split run containing a single trailing space " "
#
# Helper: replace the text of $range with $newText while forcing Word to
# keep it as its OWN run (instead of silently re-merging it into a
# neighbouring run that happens to share the same formatting). We do this
# by toggling Bold right before/after the edit -- restoring the boolean
# to its original value leaves no trace in the saved rPr, but flipping it
# momentarily is enough to make the engine treat the edited span as a
# distinct run instead of folding it back into an identically-formatted
# neighbour.
# NOTE: Word/VBA booleans read back as 0 (False) / -1 (True) -- not 1 --
# so the "flip" has to compare against 0, not 1.
function Set-RunText($range, $newText) {
    $origBold = $range.Bold
    if ($origBold -ne 0) { $range.Bold = 0 } else { $range.Bold = -1 }
    $range.Text = $newText
    $newEnd = $range.Start + $newText.Length
    $again = $word.ActiveDocument.Range($range.Start, $newEnd)
    $again.Bold = $origBold
    return $again
}

# Same trick, used just to force a run boundary at $range without changing
# its text (used to peel the trailing space off into its own <w:r>).
function Split-RunBoundary($range) {
    $origBold = $range.Bold
    if ($origBold -ne 0) { $range.Bold = 0 } else { $range.Bold = -1 }
    $range.Bold = $origBold
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "CONTRATO DE TRABAJO INDEFINIDO " title -> "...[TERMINO]"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$t1 = $p1.Range.Text
$idx1 = $t1.IndexOf("INDEFINIDO ")
$start1 = $p1.Range.Start + $idx1
$end1 = $start1 + ("INDEFINIDO ".Length)
$rng1 = $d.Range($start1, $end1)
Set-RunText $rng1 "[TERMINO]" | Out-Null

# ---------------------------------------------------------------------
# 2) "MANIZALES" / "CALDAS" place-of-work table cell
# ---------------------------------------------------------------------
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*MANIZALES*CALDAS*") {
        $p2 = $cand
        break
    }
}
if ($p2 -eq $null) { throw "Could not locate the MANIZALES / CALDAS paragraph" }

$t2 = $p2.Range.Text
$idxM = $t2.IndexOf("MANIZALES")
$startM = $p2.Range.Start + $idxM
$endM = $startM + ("MANIZALES".Length)
$rngM = $d.Range($startM, $endM)
Set-RunText $rngM "[CD_CONT]" | Out-Null

$p2b = $p2
$t2b = $p2b.Range.Text
$idxC = $t2b.IndexOf("CALDAS")
$startC = $p2b.Range.Start + $idxC
$endC = $startC + ("CALDAS".Length)
$rngC = $d.Range($startC, $endC)
Set-RunText $rngC "[DPTO_CONT]" | Out-Null

# ---------------------------------------------------------------------
# 3) "... termino indefinido y tendra vigencia ..." clause
#    "indefinido " -> "[TERMINO]" + its own trailing-space run
# ---------------------------------------------------------------------
$p3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*rmino indefinido y*") {
        $p3 = $cand
        break
    }
}
if ($p3 -eq $null) { throw "Could not locate the termino indefinido clause paragraph" }

$t3 = $p3.Range.Text
$idx3 = $t3.IndexOf("indefinido ")
$start3 = $p3.Range.Start + $idx3
$end3 = $start3 + ("indefinido ".Length)
$rng3 = $d.Range($start3, $end3)
Set-RunText $rng3 "[TERMINO] " | Out-Null

# Split the trailing space off into its own run (so it matches a freshly
# inserted <w:r> in the diff) using the same Bold-toggle trick.
$spaceStart = $start3 + "[TERMINO]".Length
$spaceEnd = $spaceStart + 1
$rngSpace = $d.Range($spaceStart, $spaceEnd)
Split-RunBoundary $rngSpace
